$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "ctb_banddis"
$ws.Range("C1").Value = "frs_banddis"
